# Payroll Suite Statutory Shared ParentalPay201718.xlsx - "Code Merge Changes" edit
# 1. Update the shared "Automation Test Result" file-path text used on H2 of the
#    weekly ShPP processing sheets (and the average-weekly-earnings report).
# 2. Clean up the stray multi-area selection (activeCellId / duplicated sqref)
#    left on the "ProcessPayrolFor16WeeklyShPP" sheet.
# 3. Move the active/selected tab from "ProcessPayrolFor48WeeklyShPP" (last sheet)
#    to "ProcessPayrolFor11WeeklyShPP" (index 5).

$wb = $excel.ActiveWorkbook

$newPath = "F:\\Automation_TestResults\\Payroll_Tax_StatutoryScenarios\\201718 Automation TestResult For Statutory Scenarios.xlsx"

$sheetsWithPath = @(
  "ProcessPayrolFor11WeeklyShPP",
  "ProcessPayrolFor12WeeklyShPP",
  "ProcessPayrolFor13WeeklyShPP",
  "ProcessPayrolFor14WeeklyShPP",
  "ProcessPayrolFor15WeeklyShPP",
  "ProcessPayrolFor16WeeklyShPP",
  "ProcessPayrolFor17WeeklyShPP",
  "ProcessPayrolFor18WeeklyShPP",
  "ProcessPayrolFor19WeeklyShPP",
  "AverageWeeklyEarningsTestReport",
  "ProcessPayrolFor46WeeklyShPP",
  "ProcessPayrolFor47WeeklyShPP",
  "ProcessPayrolFor48WeeklyShPP"
)

foreach ($name in $sheetsWithPath) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("H2").Value = $newPath
}

# Fix the odd multi-area selection (activeCell="H2" activeCellId="1" sqref="H2 H2")
# down to a plain single-cell selection on ProcessPayrolFor16WeeklyShPP.
$ws16 = $wb.Worksheets.Item("ProcessPayrolFor16WeeklyShPP")
$ws16.Activate() | Out-Null
$ws16.Range("H2").Select() | Out-Null

# Finally, move the active tab to ProcessPayrolFor11WeeklyShPP (workbook index 5),
# which also flips tabSelected off the previously active last sheet and onto this one.
$wsActive = $wb.Worksheets.Item("ProcessPayrolFor11WeeklyShPP")
$wsActive.Activate() | Out-Null
